$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 9549.549999999999
$ws.Range("I28").Value = 1693.8572
$ws.Range("J28").Value = 27879.5
$ws.Range("K28").Value = 1693.8572
$ws.Range("L28").Value = 27879.5
$ws.Range("M28").Value = -1208.8572
$ws.Range("N28").Value = -28849.5
$ws.Range("H43").Value = 1956.2142
$ws.Range("I43").Value = 1674.75
$ws.Range("J43").Value = 2068.8
$ws.Range("K43").Value = 1674.75
$ws.Range("L43").Value = 2068.8
$ws.Range("M43").Value = -1605.75
$ws.Range("N43").Value = -2206.8
$ws.Range("H53").Value = 165.57895
$ws.Range("I53").Value = 144.72728
$ws.Range("J53").Value = 194.25
$ws.Range("K53").Value = 144.72728
$ws.Range("L53").Value = 194.25
$ws.Range("M53").Value = 492.27272
$ws.Range("H62").Value = 7319.2856
$ws.Range("I62").Value = 3012.5
$ws.Range("J62").Value = 10549.375
$ws.Range("K62").Value = 3012.5
$ws.Range("L62").Value = 10549.375
$ws.Range("M62").Value = -2388.5
$ws.Range("N62").Value = -11797.375
$ws.Range("H65").Value = 7319.2856
$ws.Range("I65").Value = 3012.5
$ws.Range("J65").Value = 10549.375
$ws.Range("K65").Value = 15062.5
$ws.Range("L65").Value = 52746.875
$ws.Range("M65").Value = -11942.5
$ws.Range("N65").Value = -58986.875
$ws.Range("H76").Value = 4632591
$ws.Range("I76").Value = 5379299
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 5379299
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -5378984
$ws.Range("N76").Value = -3630
$ws.Range("H79").Value = 4632591
$ws.Range("I79").Value = 5379299
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 5379299
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -5378207
$ws.Range("N79").Value = -5184
$ws.Range("H92").Value = 8334528.5
$ws.Range("I92").Value = 10753392
$ws.Range("J92").Value = 2887.111
$ws.Range("K92").Value = 10753392
$ws.Range("L92").Value = 2887.111
$ws.Range("M92").Value = -10752144
$ws.Range("H106").Value = 8337416
$ws.Range("I106").Value = 10004249
$ws.Range("J106").Value = 3250
$ws.Range("K106").Value = 10004249
$ws.Range("L106").Value = 3250
$ws.Range("M106").Value = -10003618
$ws.Range("H135").Value = 1428.4117
$ws.Range("I135").Value = 967.4545000000001
$ws.Range("J135").Value = 2273.5
$ws.Range("K135").Value = 8707.0905
$ws.Range("L135").Value = 20461.5
$ws.Range("M135").Value = -6172.0905
$ws.Range("N135").Value = -25531.5
$ws.Range("H137").Value = 1484.2273
$ws.Range("I137").Value = 1353.5238
$ws.Range("J137").Value = 1603.5652
$ws.Range("K137").Value = 4060.5714
$ws.Range("L137").Value = 4810.6956
$ws.Range("M137").Value = -1510.5714
$ws.Range("N137").Value = -9910.695599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1203935
$ws.Range("I32").Value = 1449688.8
$ws.Range("J32").Value = 21245.062
$ws.Range("K32").Value = 1449688.8
$ws.Range("L32").Value = 21245.062
$ws.Range("M32").Value = -1449401.8
$ws.Range("N32").Value = -21819.062
$ws.Range("H45").Value = 3591
$ws.Range("I45").Value = 5000
$ws.Range("J45").Value = 2886.5
$ws.Range("K45").Value = 5000
$ws.Range("L45").Value = 2886.5
$ws.Range("M45").Value = -4623
$ws.Range("N45").Value = -3640.5
$ws.Range("H110").Value = 35218.383
$ws.Range("I110").Value = 56717.5
$ws.Range("J110").Value = 819.8
$ws.Range("K110").Value = 56717.5
$ws.Range("L110").Value = 819.8
$ws.Range("M110").Value = -54672.5
$ws.Range("N110").Value = -4909.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 18445.777
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 18445.777
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 18445.777
$ws.Range("N4").Value = -18669.777
$ws.Range("H31").Value = 4461.311
$ws.Range("I31").Value = 1116.4412
$ws.Range("J31").Value = 7304.45
$ws.Range("K31").Value = 1116.4412
$ws.Range("L31").Value = 7304.45
$ws.Range("M31").Value = -821.4412
$ws.Range("N31").Value = -7894.45
$ws.Range("H34").Value = 4461.311
$ws.Range("I34").Value = 1116.4412
$ws.Range("J34").Value = 7304.45
$ws.Range("K34").Value = 1116.4412
$ws.Range("L34").Value = 7304.45
$ws.Range("M34").Value = -914.4412
$ws.Range("N34").Value = -7708.45

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 3590
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 3590
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 10770
$ws.Range("N102").Value = -15638
$ws.Range("H105").Value = 300000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 300000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 900000
$ws.Range("N105").Value = -905242
$ws.Range("M105").ClearContents()
$ws.Range("H116").Value = 2046.3077
$ws.Range("I116").Value = 1630.8334
$ws.Range("J116").Value = 7032
$ws.Range("K116").Value = 4892.5002
$ws.Range("L116").Value = 21096
$ws.Range("M116").Value = -1450.5002
$ws.Range("N116").Value = -27980
$ws.Range("H117").Value = 1899.75
$ws.Range("I117").Value = 1249.5
$ws.Range("J117").Value = 2550
$ws.Range("K117").Value = 3748.5
$ws.Range("L117").Value = 7650
$ws.Range("M117").Value = -306.5
$ws.Range("N117").Value = -14534
$ws.Range("H118").Value = 4004.2144
$ws.Range("I118").Value = 3676.3333
$ws.Range("J118").Value = 4093.6365
$ws.Range("K118").Value = 11028.9999
$ws.Range("L118").Value = 12280.9095
$ws.Range("M118").Value = -9785.999899999999
$ws.Range("N118").Value = -14766.9095
$ws.Range("H119").Value = 1133.3334
$ws.Range("I119").Value = 1133.3334
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 3400.0002
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = 1437.9998
$ws.Range("N119").ClearContents()
$ws.Range("H120").Value = 13592.333
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 13592.333
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 40776.999
$ws.Range("N120").Value = -50452.999
$ws.Range("M120").ClearContents()
$ws.Range("H121").Value = 2021.875
$ws.Range("I121").Value = 173
$ws.Range("J121").Value = 3131.2
$ws.Range("K121").Value = 519
$ws.Range("L121").Value = 9393.599999999999
$ws.Range("M121").Value = 791
$ws.Range("N121").Value = -12013.6
$ws.Range("H123").Value = 2970
$ws.Range("I123").Value = 2970
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 8910
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -6460
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 8957.071
$ws.Range("I122").Value = 18683.166
$ws.Range("J122").Value = 1662.5
$ws.Range("K122").Value = 56049.49800000001
$ws.Range("L122").Value = 4987.5
$ws.Range("M122").Value = -53599.49800000001
$ws.Range("N122").Value = -9887.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 26176
$ws.Range("I22").Value = 875
$ws.Range("J22").Value = 36296.4
$ws.Range("K22").Value = 875
$ws.Range("L22").Value = 36296.4
$ws.Range("M22").Value = -580
$ws.Range("N22").Value = -36886.4
$ws.Range("H27").Value = 26176
$ws.Range("I27").Value = 875
$ws.Range("J27").Value = 36296.4
$ws.Range("K27").Value = 875
$ws.Range("L27").Value = 36296.4
$ws.Range("M27").Value = -768
$ws.Range("N27").Value = -36510.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6484293
$ws.Range("I132").Value = 2881.5925
$ws.Range("J132").Value = 16206410
$ws.Range("K132").Value = 8644.7775
$ws.Range("L132").Value = 48619230
$ws.Range("M132").Value = -6114.7775
$ws.Range("N132").Value = -48624290
